$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '29.679.29'
$ws.Range('E2').Value = '  -3.19%  '
Set-TextValue $ws.Range('D3') '2.096.26'
$ws.Range('E3').Value = '  -1.41%  '
$ws.Range('E4').Value = '  -0.33%  '
Set-TextValue $ws.Range('D5') '343.53'
$ws.Range('E5').Value = '  -2.58%  '
$ws.Range('E6').Value = '  -0.38%  '
Set-TextValue $ws.Range('D7') '0.5150'
$ws.Range('E7').Value = '  -2.44%  '
Set-TextValue $ws.Range('D8') '0.4406'
$ws.Range('E8').Value = '  -3.42%  '
Set-TextValue $ws.Range('D9') '53.05'
$ws.Range('E9').Value = '  -1.67%  '
Set-TextValue $ws.Range('D10') '0.09201'
$ws.Range('E10').Value = '  +1.03%  '
$ws.Range('E11').Value = '  -1.16%  '
Set-TextValue $ws.Range('D12') '24.90'
$ws.Range('E12').Value = '  +1.11%  '
Set-TextValue $ws.Range('D13') '2.093.67'
$ws.Range('E13').Value = '  -1.97%  '
Set-TextValue $ws.Range('D14') '6.750'
$ws.Range('E14').Value = '  -1.56%  '
Set-TextValue $ws.Range('D15') '8.208'
$ws.Range('E15').Value = '  +1.15%  '
Set-TextValue $ws.Range('D16') '99.39'
$ws.Range('E16').Value = '  -3.08%  '
$ws.Range('E17').Value = '  -2.47%  '
$ws.Range('E18').Value = '  -0.36%  '
Set-TextValue $ws.Range('D19') '20.69'
$ws.Range('E19').Value = '  +6.27%  '
$ws.Range('E20').Value = '  -1.18%  '
Set-TextValue $ws.Range('D21') '1.006'
$ws.Range('E21').Value = '  -0.41%  '
Set-TextValue $ws.Range('D22') '6.192'
$ws.Range('E22').Value = '  -2.64%  '
Set-TextValue $ws.Range('D23') '29.732.85'
$ws.Range('E23').Value = '  -3.23%  '
Set-TextValue $ws.Range('D24') '12.57'
$ws.Range('E24').Value = '  -2.61%  '
Set-TextValue $ws.Range('D25') '2.307'
$ws.Range('E25').Value = '  -3.40%  '
Set-TextValue $ws.Range('D26') '2.345.77'
$ws.Range('E26').Value = '  -1.68%  '
Set-TextValue $ws.Range('D27') '21.86'
$ws.Range('E27').Value = '  -2.93%  '
Set-TextValue $ws.Range('D28') '162.70'
$ws.Range('E28').Value = '  -1.18%  '
Set-TextValue $ws.Range('D29') '2.524'
$ws.Range('E29').Value = '  -1.60%  '
Set-TextValue $ws.Range('D30') '132.57'
$ws.Range('E30').Value = '  -3.26%  '
Set-TextValue $ws.Range('D31') '1.129'
$ws.Range('E31').Value = '  -5.95%  '
Set-TextValue $ws.Range('D32') '0.1048'
$ws.Range('E32').Value = '  -3.23%  '
Set-TextValue $ws.Range('D33') '1.643'
$ws.Range('E33').Value = '  -1.46%  '
Set-TextValue $ws.Range('D34') '6.156'
$ws.Range('E34').Value = '  -3.57%  '
Set-TextValue $ws.Range('D35') '3.962'
$ws.Range('E35').Value = '  -1.52%  '
Set-TextValue $ws.Range('D36') '6.051'
$ws.Range('E36').Value = '  -2.11%  '
Set-TextValue $ws.Range('D37') '10.35'
$ws.Range('E37').Value = '  +0.06%  '
Set-TextValue $ws.Range('D38') '0.02562'
$ws.Range('E38').Value = '  -3.70%  '
Set-TextValue $ws.Range('D39') '0.06711'
$ws.Range('E39').Value = '  -2.59%  '
$ws.Range('E40').Value = '  -1.39%  '
Set-TextValue $ws.Range('D41') '0.6860'
$ws.Range('E41').Value = '  -1.09%  '
Set-TextValue $ws.Range('D42') '0.2223'
$ws.Range('E42').Value = '  -4.44%  '
Set-TextValue $ws.Range('D43') '1.294'
$ws.Range('E43').Value = '  +1.14%  '
Set-TextValue $ws.Range('D44') '0.6639'
$ws.Range('E44').Value = '  +2.70%  '
Set-TextValue $ws.Range('D45') '14.14'
$ws.Range('E45').Value = '  -4.96%  '
Set-TextValue $ws.Range('D46') '2.289'
$ws.Range('E46').Value = '  -2.28%  '
Set-TextValue $ws.Range('D47') '3.613'
$ws.Range('E47').Value = '  -4.23%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D48') '0.00000000348'
$ws.Range('E48').Value = '  -5.29%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue $ws.Range('D49') '1.217'
$ws.Range('E49').Value = '  -3.23%  '
Set-TextValue $ws.Range('D50') '81.88'
$ws.Range('E50').Value = '  -1.42%  '
$ws.Range('B51').Value = 'WEMIXTOKEN'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D51') '1.160'
$ws.Range('E51').Value = '  -3.03%  '
